$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }

    $d = $ws.Cells.Item($row, 4).Value2
    $e = $ws.Cells.Item($row, 5).Value2
    $f = $ws.Cells.Item($row, 6).Value2

    $newE = $e - 1
    if ($newE -le 0) {
        $ws.Cells.Item($row, 5).Value2 = $d
        $ws.Cells.Item($row, 6).Value2 = $f + $d
    } else {
        $ws.Cells.Item($row, 5).Value2 = $newE
    }
}
